$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddProductTest")

$ws.Copy($null, $ws)
$newWs = $wb.ActiveSheet
$newWs.Name = "CartTest"

$ws.Activate()
$ws.Range("A1:C2").Select()
$newWs.Activate()
